# Apply the latest coin price / 1h-volume snapshot to Sheet1, matching the
# "Updated cryptos list" GitHub Actions commit (incl. the WEMIXToken/
# PaxDollar row-42/43 swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.941.58"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.58%  "
$c = $ws.Range("D3")
$c.Value = "'1.644.80"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.64%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$c = $ws.Range("D5")
$c.Value = "'213.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  -0.24%  "
$c = $ws.Range("D8")
$c.Value = "'23.59"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.42%  "
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("E10").Value = "  +0.27%  "
$c = $ws.Range("D11")
$c.Value = "'0.0871"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "
$c = $ws.Range("D12")
$c.Value = "'1.874.95"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'1.631.91"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E14").Value = "  +1.36%  "
$c = $ws.Range("D15")
$c.Value = "'0.566"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.56%  "
$c = $ws.Range("D16")
$c.Value = "'65.69"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "
$c = $ws.Range("D17")
$c.Value = "'27.913.28"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.43%  "
$c = $ws.Range("D18")
$c.Value = "'232.09"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "
$c = $ws.Range("D19")
$c.Value = "'7.71"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  -0.12%  "
$c = $ws.Range("D22")
$c.Value = "'10.73"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +8.15%  "
$c = $ws.Range("D23")
$c.Value = "'4.39"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("E24").Value = "  +4.08%  "
$c = $ws.Range("D25")
$c.Value = "'151.75"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$c = $ws.Range("D26")
$c.Value = "'6.94"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "
$c = $ws.Range("D27")
$c.Value = "'15.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  -0.07%  "
$c = $ws.Range("D29")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +1.14%  "
$c = $ws.Range("D31")
$c.Value = "'0.0484"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "
$c = $ws.Range("D32")
$c.Value = "'3.32"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.64%  "
$c = $ws.Range("D33")
$c.Value = "'1.459.36"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("E36").Value = "  -0.70%  "
$c = $ws.Range("D37")
$c.Value = "'0.892"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$c = $ws.Range("D38")
$c.Value = "'0.564"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.91%  "
$c = $ws.Range("D40")
$c.Value = "'0.917"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.68%  "
$c = $ws.Range("D41")
$c.Value = "'69.53"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D42")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D43")
$c.Value = "'1.01"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +0.53%  "
$c = $ws.Range("D46")
$c.Value = "'5.39"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +5.62%  "
$c = $ws.Range("D48")
$c.Value = "'1.784.11"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "
$c = $ws.Range("D49")
$c.Value = "'88.76"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("E50").Value = "  +2.21%  "
$c = $ws.Range("D51")
$c.Value = "'0.0508"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
